# Pavel - new user for linking test
# Adds a new row to the "Users" sheet for the Linking_AutoUser test account.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New row of data (row 52), mirroring the layout of the existing rows.
$ws.Range("A52").Value = "Linking_AutoUser"
$ws.Range("B52").Value = "Password1"
$ws.Range("C52").Value = ""
$ws.Range("D52").Value = ""
$ws.Range("E52").Value = "Default user for Linking tests"
$ws.Range("F52").Value = "N"
$ws.Range("G52").Value = "linking.autouser@mailinator.com"

# Match formatting (borders/style) of an existing fully-styled data row
# (row 38 uses the plain bordered style across every column, unlike rows
# that have a hyperlinked last cell).
$ws.Range("A38:G38").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the on-screen selection/scroll position to where the workbook was
# left (near the new row's sibling data, D22) like the authored view.
$null = $ws.Activate()
$null = $ws.Range("D22").Select()
$excel.ActiveWindow.ScrollRow = 13
